# Generate Report for Handback
# Updates the handback-status report: the latest handback round-trip for
# fd117f3a-506c-4a66-8548-966cf535a30d now resolves to an out-of-date
# handback file, so the zh-cn / de-de sheets gain a populated "Latest
# Target File" / "Latest Handback File" / "Latest Handback DateTime" /
# "Error Detail" for row 8, plus a widened Error Detail column.

$wb = $excel.ActiveWorkbook

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b5b9c95ee1dac8e137049e529776cfa8ae14a0e/e2e/fd117f3a-506c-4a66-8548-966cf535a30d.md"
$targetDisplay = "fd117f3a-506c-4a66-8548-966cf535a30d.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/45207ae6b0afbe8f54dd743943057e6dcf93baf3/e2e/fd117f3a-506c-4a66-8548-966cf535a30d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b5b9c95ee1dac8e137049e529776cfa8ae14a0e/e2e/fd117f3a-506c-4a66-8548-966cf535a30d.md."

# Widened "Error Detail" column (col 16 / P) on the per-language sheets.
$newErrorColWidth = 39.166666666666664

function Update-LanguageSheet($sheetName, $handbackXlf, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column P ("Error Detail") widened to fit the long message.
    $ws.Columns.Item(16).ColumnWidth = $newErrorColWidth

    # Row 8 ("fd117f3a-506c-4a66-8548-966cf535a30d.md") now has a target
    # file, handback file, handback datetime, and error detail populated.
    $i8 = $ws.Range("I8")
    $ws.Hyperlinks.Add($i8, $hyperlinkTarget, "", "", $targetDisplay)
    $i8.Font.Underline = 2
    $i8.Font.Color = 15570276

    $ws.Range("J8").Value = $handbackXlf
    $ws.Range("K8").Value = $handbackDateTime
    $ws.Range("P8").Value = $errorDetail
}

Update-LanguageSheet "zh-cn" "fd117f3a-506c-4a66-8548-966cf535a30d.48b52e17847ebe53b67608aae9fb3a87e42aad5e.zh-cn.xlf" "2016-08-19 22:49:24"
Update-LanguageSheet "de-de" "fd117f3a-506c-4a66-8548-966cf535a30d.48b52e17847ebe53b67608aae9fb3a87e42aad5e.de-de.xlf" "2016-08-19 22:49:31"
